# Adding Latin names for pg. 180
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Latin Name) values for rows 43-64. These must be written in
# the same order the source workbook's shared-string table was built in
# (new species names first appear row-by-row in a particular non-linear
# sequence), so that the resulting xl/sharedStrings.xml unique-string
# ordering matches exactly.
$latinRows = @(49,64,43,47,48,51,55,57,60,61,62,63,44,45,46,50,52,54,53,56,58,59)
$latinNames = @("NA","Microcerculus marginatus","Phaethornis superciliosus","Pipra fasciicauda","Turdus fumigatus","Basileuterus rivularis","Phaethornis ruber","Chloroceryle inda","Xiphorphynchus guttatus","Philydor pyrrhodes","Thamnophilus aethiops","Pipra iris","Notharctus tectus","Pteroglossus aracari","Cotinga cayana","Cyanerpes cyaneus","Euphonia violacea","Nyctidromus albicollis","Piaya minuta","Trogon violaceus","Pteroglossus bitorquatus","Veniliornis affinis")

for ($i = 0; $i -lt $latinRows.Length; $i++) {
    $ws.Cells.Item($latinRows[$i], 2).Value = $latinNames[$i]
}

# Row 64 was missing its frequency-of-capture data (columns D:CB) in the
# source workbook; fill it in now, matching row 63's pattern, and add the
# row total in CC64.
$row64Data = @(0,1,2,3,1,1,5,0,3,0,0,1,5,2,4,0,4,4,1,5,2,3,1,3,1,3,2,0,0,0,0,0,0,1,3,4,0,3,1,1,0,3,1,0,0,0,0,0,0,0,0,0,0,0,0,0,1,2,5,4,1,0,1,9,0,0,0,0,0,0,1,0,0,0,0,1,0)

for ($i = 0; $i -lt $row64Data.Length; $i++) {
    $ws.Cells.Item(64, 4 + $i).Value = $row64Data[$i]
}

$ws.Range("CC64").Formula = "=SUM(D64:CB64)"

# Restore the view: unfreeze scroll back to the top and select C61.
$ws.Activate()
$ws.Range("C61").Select()
